$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2288135593220339
$ws.Range("C2").Value = 0.5169491525423728
$ws.Range("J2").Value = 0.01271186440677966
$ws.Range("P2").Value = 0.1440677966101695
$ws.Range("S2").Value = 0.09745762711864407
$ws.Range("C3").Value = 0.02419354838709677
$ws.Range("J3").Value = 0.02419354838709677
$ws.Range("P3").Value = 0.7741935483870968
$ws.Range("S3").Value = 0.1774193548387097
$ws.Range("J4").Value = 0.02941176470588235
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2941176470588235
$ws.Range("J5").Value = 0.4
$ws.Range("P5").Value = 0.6
$ws.Range("B6").Value = 0.04184100418410042
$ws.Range("D6").Value = 0.01255230125523013
$ws.Range("F6").Value = 0.08368200836820083
$ws.Range("J6").Value = 0.2301255230125523
$ws.Range("O6").Value = 0.01673640167364017
$ws.Range("Q6").Value = 0.1548117154811715
$ws.Range("R6").Value = 0.07949790794979079
$ws.Range("S6").Value = 0.3807531380753138
$ws.Range("B7").Value = 0.1067415730337079
$ws.Range("D7").Value = 0.02247191011235955
$ws.Range("F7").Value = 0.0898876404494382
$ws.Range("J7").Value = 0.0898876404494382
$ws.Range("O7").Value = 0.01685393258426966
$ws.Range("Q7").Value = 0.1348314606741573
$ws.Range("R7").Value = 0.1348314606741573
$ws.Range("S7").Value = 0.4044943820224719
$ws.Range("B8").Value = 0.0860655737704918
$ws.Range("D8").Value = 0.00819672131147541
$ws.Range("E8").Value = 0.006147540983606557
$ws.Range("F8").Value = 0.05327868852459016
$ws.Range("J8").Value = 0.1086065573770492
$ws.Range("O8").Value = 0.00819672131147541
$ws.Range("Q8").Value = 0.1741803278688525
$ws.Range("R8").Value = 0.1045081967213115
$ws.Range("S8").Value = 0.4508196721311475
$ws.Range("B9").Value = 0.05
$ws.Range("D9").Value = 0.03333333333333333
$ws.Range("F9").Value = 0.06111111111111111
$ws.Range("J9").Value = 0.09444444444444444
$ws.Range("O9").Value = 0.02777777777777778
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.09444444444444444
$ws.Range("S9").Value = 0.4388888888888889
$ws.Range("B10").Value = 0.09083044982698962
$ws.Range("D10").Value = 0.01557093425605536
$ws.Range("E10").Value = 0.00259515570934256
$ws.Range("F10").Value = 0.08564013840830449
$ws.Range("J10").Value = 0.09775086505190311
$ws.Range("O10").Value = 0.0198961937716263
$ws.Range("Q10").Value = 0.2136678200692042
$ws.Range("R10").Value = 0.09688581314878893
$ws.Range("S10").Value = 0.3771626297577855
$ws.Range("G11").Value = 0.1319444444444444
$ws.Range("J11").Value = 0.1284722222222222
$ws.Range("K11").Value = 0.2048611111111111
$ws.Range("L11").Value = 0.5173611111111112
$ws.Range("S11").Value = 0.01736111111111111
$ws.Range("G12").Value = 0.7124183006535948
$ws.Range("J12").Value = 0.2418300653594771
$ws.Range("K12").Value = 0.0130718954248366
$ws.Range("L12").Value = 0.0261437908496732
$ws.Range("S12").Value = 0.006535947712418301
$ws.Range("G13").Value = 0.7948717948717948
$ws.Range("J13").Value = 0.2051282051282051
$ws.Range("F15").Value = 0.02604166666666667
$ws.Range("H15").Value = 0.2135416666666667
$ws.Range("I15").Value = 0.08333333333333333
$ws.Range("J15").Value = 0.3385416666666667
$ws.Range("K15").Value = 0.05208333333333334
$ws.Range("M15").Value = 0.01041666666666667
$ws.Range("O15").Value = 0.02604166666666667
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.02649006622516556
$ws.Range("H16").Value = 0.152317880794702
$ws.Range("I16").Value = 0.08609271523178808
$ws.Range("J16").Value = 0.4105960264900662
$ws.Range("K16").Value = 0.1258278145695364
$ws.Range("M16").Value = 0.01986754966887417
$ws.Range("N16").Value = 0.006622516556291391
$ws.Range("O16").Value = 0.04635761589403974
$ws.Range("S16").Value = 0.1258278145695364
$ws.Range("F17").Value = 0.01411764705882353
$ws.Range("H17").Value = 0.1811764705882353
$ws.Range("I17").Value = 0.08470588235294117
$ws.Range("J17").Value = 0.4494117647058823
$ws.Range("K17").Value = 0.09411764705882353
$ws.Range("M17").Value = 0.01176470588235294
$ws.Range("N17").Value = 0.004705882352941176
$ws.Range("O17").Value = 0.04470588235294118
$ws.Range("S17").Value = 0.1152941176470588
$ws.Range("F18").Value = 0.03139013452914798
$ws.Range("H18").Value = 0.2331838565022422
$ws.Range("I18").Value = 0.04932735426008968
$ws.Range("J18").Value = 0.42152466367713
$ws.Range("K18").Value = 0.08520179372197309
$ws.Range("M18").Value = 0.0179372197309417
$ws.Range("N18").Value = 0.004484304932735426
$ws.Range("O18").Value = 0.06278026905829596
$ws.Range("S18").Value = 0.09417040358744394
$ws.Range("F19").Value = 0.01631321370309951
$ws.Range("H19").Value = 0.2398042414355628
$ws.Range("I19").Value = 0.08564437194127243
$ws.Range("J19").Value = 0.3466557911908646
$ws.Range("K19").Value = 0.1052202283849918
$ws.Range("M19").Value = 0.02120717781402937
$ws.Range("N19").Value = 0.0008156606851549756
$ws.Range("O19").Value = 0.06525285481239804
$ws.Range("S19").Value = 0.1190864600326264
